$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Header row: "sig(0.9)" -> "sig(0.90)" for E1:J1 ---
# This also drops "sig(0.9)" from the shared-strings table (no longer referenced
# anywhere), which causes Sheet2!Q1/R1 and Sheet4!A1 to be re-indexed
# automatically on save.
$ws3.Range("E1:J1").Value = "sig(0.90)"

# --- Body grid: move each row's 0.5 marker to its new column ---
# Rows 2-7 already have the right cell populated; only rows 8-19 change.
$ws3.Range("G8").Value = 0.5

$ws3.Range("G9").ClearContents()
$ws3.Range("F9").Value = 0.5

$ws3.Range("F10").ClearContents()
$ws3.Range("G10").Value = 0.5

$ws3.Range("G11").ClearContents()
$ws3.Range("H11").Value = 0.5

$ws3.Range("H12").ClearContents()
$ws3.Range("F12").Value = 0.5

$ws3.Range("F13").ClearContents()
$ws3.Range("G13").Value = 0.5

$ws3.Range("G14").ClearContents()
$ws3.Range("H14").Value = 0.5

$ws3.Range("H15").ClearContents()
$ws3.Range("F15").Value = 0.5

$ws3.Range("F16").ClearContents()
$ws3.Range("G16").Value = 0.5

$ws3.Range("G17").ClearContents()
$ws3.Range("H17").Value = 0.5

$ws3.Range("H18").ClearContents()
$ws3.Range("I18").Value = 0.5

$ws3.Range("I19").ClearContents()
$ws3.Range("J19").Value = 0.5

# --- Row 20 (old J20 = 0.5) is dropped entirely ---
$ws3.Rows.Item(20).Delete()

# --- Selection moves from E1 to J20 ---
$ws3.Activate()
$ws3.Range("J20").Select()
